$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row 2
$ws.Range("C2").Value = 1110
$ws.Range("D2").Value = "d78328s"
$ws.Range("E2").Value = 879

# Date cells: set format on F2 first, then copy its format onto G2 so both
# share a single cellXfs style entry (numFmtId 14, i.e. built-in short date).
$ws.Range("F2").Value = 43628
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Value = 43636
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("H2").Value = 109230
$ws.Range("I2").Value = 888
$ws.Range("J2").Value = 8239
$ws.Range("K2").Value = "bison"

# Column width adjustments for columns I and J (values chosen so the
# engine's internal pixel-quantized width lands on the closest achievable
# approximation of the target widths 22.85546875 / 5.140625)
$ws.Columns.Item(9).ColumnWidth = 22
$ws.Columns.Item(10).ColumnWidth = 4.333

# Move selection to the newly added cell
$ws.Range("K2").Select() | Out-Null
